$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 71, shifting existing rows 71..104 down to 72..105
$ws.Rows.Item(71).Insert()

# Fill in the new row 71 with the new record's values
$ws.Cells.Item(71, 1).Value = 5
$ws.Cells.Item(71, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(71, 3).Value = "Maule"
$ws.Cells.Item(71, 4).Value = (Get-Date -Year 2023 -Month 7 -Day 27 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(71, 5).Value = 7
$ws.Cells.Item(71, 6).Value = 100112040
$ws.Cells.Item(71, 7).Value = "Cilantro"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 150
$ws.Cells.Item(71, 11).Value = 9000
$ws.Cells.Item(71, 12).Value = 9000
$ws.Cells.Item(71, 13).Value = 9000
$ws.Cells.Item(71, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(71, 15).Value = "Región Metropolitana"
$ws.Cells.Item(71, 16).Value = 250
$ws.Cells.Item(71, 17).Value = 36
$ws.Cells.Item(71, 18).Value = "Hortaliza"
